# Corrected excel sheets for application fix issues
#
# Re-applies (via Excel COM interop) the view/selection state and the two
# data corrections captured by the source diff:
#   - "Repayment schedule"!K2 / L2 : 0 -> 100
#   - "Transactions"!A2 / A3       : 112 -> 387, 111 -> 386
#   - updated cell selections on NewLoanInput, Summary, Repayment schedule
#     and Transactions, with "NewLoanInput" ending up as the active/visible
#     sheet (previously "Repayment schedule" was active).

$wb = $excel.ActiveWorkbook

# --- Data value changes -----------------------------------------------

$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("K2").Value = 100
$wsRepay.Range("L2").Value = 100

$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 387
$wsTrans.Range("A3").Value = 386

# --- View / selection changes -------------------------------------------
# (Writing cell values above does not disturb any sheet's view/selection
# state, so the order of the two phases does not matter.)

# Summary: selection -> E4
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate() | Out-Null
$wsSummary.Range("E4").Select() | Out-Null

# Repayment schedule: selection -> L2 (also no longer the active tab)
$wsRepay.Activate() | Out-Null
$wsRepay.Range("L2").Select() | Out-Null

# Transactions: selection -> A2:XFD4
$wsTrans.Activate() | Out-Null
$wsTrans.Range("A2:XFD4").Select() | Out-Null

# NewLoanInput: becomes the active/visible sheet, selection -> B19.
# Activating this sheet last makes it the workbook's active tab.
$wsNewLoan = $wb.Worksheets.Item("NewLoanInput")
$wsNewLoan.Activate() | Out-Null
$wsNewLoan.Range("B19").Select() | Out-Null
